$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Consumption (A) and Timestamp (B) values for rows 2-27 (new Dabaca location data)
$ws.Range("A2").Value = 5283
$ws.Range("B2").Value = 45807
$ws.Range("A3").Value = 5275
$ws.Range("B3").Value = 45807.01041666666
$ws.Range("A4").Value = 5255
$ws.Range("B4").Value = 45807.02083333334
$ws.Range("A5").Value = 5187
$ws.Range("B5").Value = 45807.03125
$ws.Range("A6").Value = 5103
$ws.Range("B6").Value = 45807.04166666666
$ws.Range("A7").Value = 5067
$ws.Range("B7").Value = 45807.05208333334
$ws.Range("A8").Value = 5077
$ws.Range("B8").Value = 45807.0625
$ws.Range("A9").Value = 5042
$ws.Range("B9").Value = 45807.07291666666
$ws.Range("A10").Value = 5058
$ws.Range("B10").Value = 45807.08333333334
$ws.Range("A11").Value = 5026
$ws.Range("B11").Value = 45807.09375
$ws.Range("A12").Value = 5113
$ws.Range("B12").Value = 45807.10416666666
$ws.Range("A13").Value = 5056
$ws.Range("B13").Value = 45807.11458333334
$ws.Range("A14").Value = 5007
$ws.Range("B14").Value = 45807.125
$ws.Range("A15").Value = 5017
$ws.Range("B15").Value = 45807.13541666666
$ws.Range("A16").Value = 5048
$ws.Range("B16").Value = 45807.14583333334
$ws.Range("A17").Value = 5010
$ws.Range("B17").Value = 45807.15625
$ws.Range("A18").Value = 5047
$ws.Range("B18").Value = 45807.16666666666
$ws.Range("A19").Value = 5045
$ws.Range("B19").Value = 45807.17708333334
$ws.Range("A20").Value = 5055
$ws.Range("B20").Value = 45807.1875
$ws.Range("A21").Value = 5068
$ws.Range("B21").Value = 45807.19791666666
$ws.Range("A22").Value = 5286
$ws.Range("B22").Value = 45807.20833333334
$ws.Range("A23").Value = 5405
$ws.Range("B23").Value = 45807.21875
$ws.Range("A24").Value = 5477
$ws.Range("B24").Value = 45807.22916666666
$ws.Range("A25").Value = 5601
$ws.Range("B25").Value = 45807.23958333334
$ws.Range("A26").Value = 5852
$ws.Range("B26").Value = 45807.25
$ws.Range("A27").Value = 6007
$ws.Range("B27").Value = 45807.26041666666

# Remove the now-unused trailing rows 28-44
$ws.Range("A28:B44").EntireRow.Delete()
